$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    3 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4 = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 13.86384647080068, 14.74297629124414)
    5 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    6 = @(0.1169995834814548, 0.04103571897497393, 3.223369029078222, 0.5333859586016987, 3.914790290136349)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
